$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text while assigning values,
# so numeric-looking strings like "231.77" are not auto-converted to numbers.
# Column E values always contain a % sign so they remain text naturally.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.711.31"
$ws.Range("E2").Value = "  -1.65%  "
$ws.Range("D3").Value = "1.803.98"
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "231.77"
$ws.Range("E5").Value = "  -1.71%  "
$ws.Range("D6").Value = "0.5948"
$ws.Range("E6").Value = "  -2.51%  "
$ws.Range("D8").Value = "0.2787"
$ws.Range("E8").Value = "  -0.55%  "
$ws.Range("D9").Value = "0.06835"
$ws.Range("D10").Value = "23.36"
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("D11").Value = "0.07535"
$ws.Range("E11").Value = "  -1.62%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "4.784"
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.718.48"
$ws.Range("E13").Value = "  -6.28%  "
$ws.Range("D14").Value = "0.6239"
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("D15").Value = "2.048.67"
$ws.Range("E15").Value = "  -1.17%  "
$ws.Range("D16").Value = "0.000009325"
$ws.Range("E16").Value = "  -6.72%  "
$ws.Range("D17").Value = "75.61"
$ws.Range("E17").Value = "  -3.79%  "
$ws.Range("D18").Value = "28.665.57"
$ws.Range("D19").Value = "5.496"
$ws.Range("E19").Value = "  -6.20%  "
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").Value = "210.22"
$ws.Range("E21").Value = "  -7.19%  "
$ws.Range("D22").Value = "11.47"
$ws.Range("E22").Value = "  -2.55%  "
$ws.Range("D23").Value = "6.865"
$ws.Range("E23").Value = "  -1.93%  "
$ws.Range("D24").Value = "1.004"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").Value = "154.29"
$ws.Range("E25").Value = "  -0.77%  "
$ws.Range("D26").Value = "7.875"
$ws.Range("E26").Value = "  -2.14%  "
$ws.Range("D27").Value = "0.1274"
$ws.Range("E27").Value = "  -2.40%  "
$ws.Range("D28").Value = "16.39"
$ws.Range("E28").Value = "  -1.08%  "
$ws.Range("D29").Value = "1.427"
$ws.Range("E29").Value = "  -4.55%  "
$ws.Range("D30").Value = "0.06192"
$ws.Range("E30").Value = "  -2.52%  "
$ws.Range("E31").Value = "  -2.08%  "
$ws.Range("D32").Value = "3.784"
$ws.Range("E32").Value = "  -0.90%  "
$ws.Range("D33").Value = "3.753"
$ws.Range("E33").Value = "  -1.23%  "
$ws.Range("D34").Value = "1.722"
$ws.Range("E34").Value = "  -0.82%  "
$ws.Range("D35").Value = "1.066"
$ws.Range("E35").Value = "  -5.14%  "
$ws.Range("D36").Value = "0.6405"
$ws.Range("E36").Value = "  -0.82%  "
$ws.Range("D37").Value = "2.493"
$ws.Range("E37").Value = "  -1.96%  "
$ws.Range("D38").Value = "2.714"
$ws.Range("E38").Value = "  -0.44%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "6.476"
$ws.Range("E39").Value = "  -1.12%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.01714"
$ws.Range("E40").Value = "  -1.33%  "
$ws.Range("D41").Value = "1.134.61"
$ws.Range("E41").Value = "  -6.57%  "
$ws.Range("D42").Value = "0.8771"
$ws.Range("E42").Value = "  -3.84%  "
$ws.Range("D43").Value = "1.009"
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("D44").Value = "100.75"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").Value = "1.965.54"
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("E47").Value = "  -5.18%  "
$ws.Range("E48").Value = "  -0.35%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "8.348"
$ws.Range("E49").Value = "  -2.58%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.05475"
$ws.Range("E50").Value = "  -0.80%  "
$ws.Range("D51").Value = "0.4485"
$ws.Range("E51").Value = "  -1.72%  "

# Restore the original (default/General) style so no visible formatting changes remain.
$priceRange.Style = "Normal"
